# Update the lattice multiplication exercise values to match the
# author's c986bee output. Each table cell holds a run with five
# <w:t> text fragments joined by <w:br/> line breaks:
#   "AA x BB" / "  C    D" / "  ----" / "E|    |" / "F|    |"
# Word's Range.Text exposes each <w:br/> as a vertical-tab char
# ([char]11), so we rebuild each cell's Range.Text in one assignment
# to preserve the run's sz=32 formatting and the break structure.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellText($table, $row, $col, $newText) {
    $cell = $table.Cell($row, $col)
    $cell.Range.Text = $newText
}

$txt = "50 x 16" + [char]11 + "  1    6" + [char]11 + "  ----" + [char]11 + "5|    |" + [char]11 + "0|    |"
Set-CellText $t 1 1 $txt
$txt = "50 x 17" + [char]11 + "  1    7" + [char]11 + "  ----" + [char]11 + "5|    |" + [char]11 + "0|    |"
Set-CellText $t 1 2 $txt
$txt = "97 x 57" + [char]11 + "  5    7" + [char]11 + "  ----" + [char]11 + "9|    |" + [char]11 + "7|    |"
Set-CellText $t 1 3 $txt
$txt = "35 x 88" + [char]11 + "  8    8" + [char]11 + "  ----" + [char]11 + "3|    |" + [char]11 + "5|    |"
Set-CellText $t 2 1 $txt
$txt = "36 x 95" + [char]11 + "  9    5" + [char]11 + "  ----" + [char]11 + "3|    |" + [char]11 + "6|    |"
Set-CellText $t 2 2 $txt
$txt = "94 x 61" + [char]11 + "  6    1" + [char]11 + "  ----" + [char]11 + "9|    |" + [char]11 + "4|    |"
Set-CellText $t 2 3 $txt
$txt = "19 x 85" + [char]11 + "  8    5" + [char]11 + "  ----" + [char]11 + "1|    |" + [char]11 + "9|    |"
Set-CellText $t 3 1 $txt
$txt = "55 x 90" + [char]11 + "  9    0" + [char]11 + "  ----" + [char]11 + "5|    |" + [char]11 + "5|    |"
Set-CellText $t 3 2 $txt
$txt = "76 x 77" + [char]11 + "  7    7" + [char]11 + "  ----" + [char]11 + "7|    |" + [char]11 + "6|    |"
Set-CellText $t 3 3 $txt
$txt = "82 x 19" + [char]11 + "  1    9" + [char]11 + "  ----" + [char]11 + "8|    |" + [char]11 + "2|    |"
Set-CellText $t 4 1 $txt
$txt = "23 x 43" + [char]11 + "  4    3" + [char]11 + "  ----" + [char]11 + "2|    |" + [char]11 + "3|    |"
Set-CellText $t 4 2 $txt
$txt = "34 x 99" + [char]11 + "  9    9" + [char]11 + "  ----" + [char]11 + "3|    |" + [char]11 + "4|    |"
Set-CellText $t 4 3 $txt
$txt = "21 x 58" + [char]11 + "  5    8" + [char]11 + "  ----" + [char]11 + "2|    |" + [char]11 + "1|    |"
Set-CellText $t 5 1 $txt
$txt = "24 x 79" + [char]11 + "  7    9" + [char]11 + "  ----" + [char]11 + "2|    |" + [char]11 + "4|    |"
Set-CellText $t 5 2 $txt
$txt = "72 x 82" + [char]11 + "  8    2" + [char]11 + "  ----" + [char]11 + "7|    |" + [char]11 + "2|    |"
Set-CellText $t 5 3 $txt
